$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text,
# matching the original inline-string cell type.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.177.21'
$ws.Range('E2').Value = '  +9.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.680.66'
$ws.Range('E3').Value = '  +5.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.92'
$ws.Range('E5').Value = '  +7.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9977'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3722'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3452'
$ws.Range('E8').Value = '  +3.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.25'
$ws.Range('E9').Value = '  +10.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.196'
$ws.Range('E10').Value = '  +3.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07300'
$ws.Range('E11').Value = '  +3.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9989'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.47'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.134'
$ws.Range('E14').Value = '  +3.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.782'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.678.28'
$ws.Range('E16').Value = '  +6.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001110'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9980'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06717'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '81.65'
$ws.Range('E20').Value = '  +6.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.57'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.134'
$ws.Range('E22').Value = '  +3.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.02'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.148.87'
$ws.Range('E24').Value = '  +9.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.413'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.364'
$ws.Range('E26').Value = '  -9.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.687'
$ws.Range('E27').Value = '  +7.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.74'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.61'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.863.76'
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.12'
$ws.Range('E31').Value = '  +5.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.362'
$ws.Range('E32').Value = '  +8.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.036'
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9848'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.761'
$ws.Range('E35').Value = '  +6.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08467'
$ws.Range('E36').Value = '  +2.60%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.072'
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '12.32'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06446'
$ws.Range('E39').Value = '  +3.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.378'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02353'
$ws.Range('E41').Value = '  +6.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.264'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2133'
$ws.Range('E43').Value = '  +5.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6179'
$ws.Range('E44').Value = '  +2.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9987'
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.26'
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.811'
$ws.Range('E47').Value = '  +4.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5974'
$ws.Range('E48').Value = '  +3.52%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '128.01'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.047'
$ws.Range('E50').Value = '  +5.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07203'
$ws.Range('E51').Value = '  +5.52%  '
